# Atualizacao de tarefas a serem cumpridas
#
# Appends new ", <tarefa>" fragments to the end of five task-assignment
# paragraphs (one per team member) at the top of the document, mirroring
# the existing "Projetos atualizado no GitHub / Documentacao do Projeto
# Atualizada" run with a comma + (usually bold) task name.

$d = $word.ActiveDocument

function Append-Runs {
    param(
        [int]$ParagraphIndex,
        [object[]]$Pieces   # array of @{ Text = "..."; Bold = $true/$false }
    )

    Write-Host "Append-Runs called with ParagraphIndex=$ParagraphIndex Pieces=$($Pieces.Count)"
    $para = $d.Paragraphs.Item($ParagraphIndex)
    $range = $para.Range
    # Exclude the trailing paragraph mark so new runs land inside the
    # paragraph, right after its current last run.
    $range.End = $range.End - 1

    foreach ($piece in $Pieces) {
        $insertPos = $range.End
        $range.InsertAfter($piece.Text)
        if ($piece.Bold) {
            $newRange = $d.Range($insertPos, $insertPos + $piece.Text.Length)
            $newRange.Font.Bold = $true
        }
        # Keep $range collapsed/growing at the paragraph's (new) end so the
        # next piece is appended after this one.
        $range.End = $insertPos + $piece.Text.Length
    }
}

# Vitor Xavier -> ", " + "Modelagem Logica do Projeto v1" (bold), with the
# leading "," plain and the following space bold (matches source diff).
Append-Runs 3 @(
    @{ Text = ","; Bold = $false },
    @{ Text = " "; Bold = $true },
    @{ Text = "Modelagem Lógica do Projeto v1"; Bold = $true }
)

# Mateus Resende -> ", " + "Site Estatico Institucional ... JavaScript   " (bold)
Append-Runs 4 @(
    @{ Text = ", "; Bold = $false },
    @{ Text = "Site Estático Institucional – Local em HTML/CSS/"; Bold = $true },
    @{ Text = "JavaScript"; Bold = $true },
    @{ Text = "   "; Bold = $true }
)

# Vitoria Mirella -> "," + " " (both plain, separate runs) + bold site text
Append-Runs 5 @(
    @{ Text = ","; Bold = $false },
    @{ Text = " "; Bold = $false },
    @{ Text = "Site Estático Institucional – Local em HTML/CSS/"; Bold = $true },
    @{ Text = "JavaScript"; Bold = $true },
    @{ Text = "   "; Bold = $true }
)

# Thamiris Ayumi -> ", " + "Revisar a calculadora" (bold)
Append-Runs 6 @(
    @{ Text = ", "; Bold = $false },
    @{ Text = "Revisar a calculadora"; Bold = $true }
)

# Victor Augusto -> ", " + "Modelagem Logica do Projeto v1" (bold)
Append-Runs 7 @(
    @{ Text = ", "; Bold = $false },
    @{ Text = "Modelagem Lógica do Projeto v1"; Bold = $true }
)
